$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 04:52"

# Row 45: Australia
$ws.Cells.Item(45, 1).Value = "Australia"
$ws.Cells.Item(45, 2).Value = 6660
$ws.Cells.Item(45, 3).Value = 11
$ws.Cells.Item(45, 4).Value = 5034
$ws.Cells.Item(45, 5).Value = 1551
$ws.Cells.Item(45, 6).Value = 46
$ws.Cells.Item(45, 7).Value = 1
$ws.Cells.Item(45, 8).Value = 75

# Row 74: Nueva Zelanda
$ws.Cells.Item(74, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(74, 2).Value = 1451
$ws.Cells.Item(74, 3).Value = 3
$ws.Cells.Item(74, 4).Value = 1065
$ws.Cells.Item(74, 5).Value = 370
$ws.Cells.Item(74, 6).Value = 2
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = 16

# Row 97: Bolivia
$ws.Cells.Item(97, 1).Value = "Bolivia"
$ws.Cells.Item(97, 2).Value = 672
$ws.Cells.Item(97, 3).Value = 63
$ws.Cells.Item(97, 4).Value = 44
$ws.Cells.Item(97, 5).Value = 588
$ws.Cells.Item(97, 6).Value = 3
$ws.Cells.Item(97, 7).Value = 3
$ws.Cells.Item(97, 8).Value = 40

# Row 98: Niger
$ws.Cells.Item(98, 1).Value = "Niger"
$ws.Cells.Item(98, 2).Value = 662
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 193
$ws.Cells.Item(98, 5).Value = 447
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 22

# Row 99: Albania
$ws.Cells.Item(99, 1).Value = "Albania"
$ws.Cells.Item(99, 2).Value = 634
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 356
$ws.Cells.Item(99, 5).Value = 251
$ws.Cells.Item(99, 6).Value = 4
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 27

# Row 100: Kirguistan
$ws.Cells.Item(100, 1).Value = "Kirguistan"
$ws.Cells.Item(100, 2).Value = 612
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 254
$ws.Cells.Item(100, 5).Value = 351
$ws.Cells.Item(100, 6).Value = 5
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 7

# Row 103: Honduras
$ws.Cells.Item(103, 1).Value = "Honduras"
$ws.Cells.Item(103, 2).Value = 519
$ws.Cells.Item(103, 3).Value = 9
$ws.Cells.Item(103, 4).Value = 31
$ws.Cells.Item(103, 5).Value = 441
$ws.Cells.Item(103, 6).Value = 10
$ws.Cells.Item(103, 7).Value = 1
$ws.Cells.Item(103, 8).Value = 47

# Row 113: Guatemala
$ws.Cells.Item(113, 1).Value = "Guatemala"
$ws.Cells.Item(113, 2).Value = 342
$ws.Cells.Item(113, 3).Value = 26
$ws.Cells.Item(113, 4).Value = 25
$ws.Cells.Item(113, 5).Value = 307
$ws.Cells.Item(113, 6).Value = 3
$ws.Cells.Item(113, 7).Value = 2
$ws.Cells.Item(113, 8).Value = 10

# Row 114: Sri Lanka
$ws.Cells.Item(114, 1).Value = "Sri Lanka"
$ws.Cells.Item(114, 2).Value = 330
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 105
$ws.Cells.Item(114, 5).Value = 218
$ws.Cells.Item(114, 6).Value = 2
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 7

# Row 115: Mauricio
$ws.Cells.Item(115, 1).Value = "Mauricio"
$ws.Cells.Item(115, 2).Value = 329
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 261
$ws.Cells.Item(115, 5).Value = 59
$ws.Cells.Item(115, 6).Value = 3
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 9

# Row 116: Mayotte
$ws.Cells.Item(116, 1).Value = "Mayotte"
$ws.Cells.Item(116, 2).Value = 326
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 125
$ws.Cells.Item(116, 5).Value = 197
$ws.Cells.Item(116, 6).Value = 4
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 4

# Row 126: El Salvador
$ws.Cells.Item(126, 1).Value = "El Salvador"
$ws.Cells.Item(126, 2).Value = 237
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(126, 4).Value = 63
$ws.Cells.Item(126, 5).Value = 166
$ws.Cells.Item(126, 6).Value = 3
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 8

# Row 132: Sudan
$ws.Cells.Item(132, 1).Value = "Sudan"
$ws.Cells.Item(132, 2).Value = 162
$ws.Cells.Item(132, 3).Value = 22
$ws.Cells.Item(132, 4).Value = 14
$ws.Cells.Item(132, 5).Value = 135
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 13

# Row 137: Birmania
$ws.Cells.Item(137, 1).Value = "Birmania"
$ws.Cells.Item(137, 2).Value = 127
$ws.Cells.Item(137, 3).Value = 4
$ws.Cells.Item(137, 4).Value = 7
$ws.Cells.Item(137, 5).Value = 115
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 5

# Row 154: San Martin (Parte Holandesa)
$ws.Cells.Item(154, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(154, 2).Value = 73
$ws.Cells.Item(154, 3).Value = 2
$ws.Cells.Item(154, 4).Value = 22
$ws.Cells.Item(154, 5).Value = 39
$ws.Cells.Item(154, 6).Value = 8
$ws.Cells.Item(154, 7).Value = 1
$ws.Cells.Item(154, 8).Value = 12
